# Auto-generated Excel COM-interop script
# Updates raw market-price data cells across multiple worksheets
# to match the refreshed values pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 3339000
$ws.Range("I20").Value = 3339000
$ws.Range("K20").Value = 3339000
$ws.Range("M20").Value = -3338770
$ws.Range("H32").Value = 3177.9
$ws.Range("J32").Value = 3253.2222
$ws.Range("L32").Value = 3253.2222
$ws.Range("N32").Value = -3905.2222
$ws.Range("H35").Value = 3339000
$ws.Range("I35").Value = 3339000
$ws.Range("K35").Value = 3339000
$ws.Range("M35").Value = -3338621
$ws.Range("H116").Value = 19136116
$ws.Range("I116").Value = 24354188
$ws.Range("J116").Value = 3181.6667
$ws.Range("K116").Value = 24354188
$ws.Range("L116").Value = 3181.6667
$ws.Range("M116").Value = -24350746
$ws.Range("N116").Value = -10065.6667
$ws.Range("H132").Value = 5194.9756
$ws.Range("I132").Value = 1466.6945
$ws.Range("K132").Value = 4400.083500000001
$ws.Range("M132").Value = -1870.083500000001
$ws.Range("H137").Value = 185076
$ws.Range("I137").Value = 297101.47
$ws.Range("J137").Value = 3701.4285
$ws.Range("K137").Value = 891304.4099999999
$ws.Range("L137").Value = 11104.2855
$ws.Range("M137").Value = -888754.4099999999
$ws.Range("N137").Value = -16204.2855
$ws.Range("H138").Value = 7170.411
$ws.Range("J138").Value = 7988.364
$ws.Range("L138").Value = 23965.092
$ws.Range("N138").Value = -34245.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 700.9375
$ws.Range("I2").Value = 621.8276
$ws.Range("K2").Value = 621.8276
$ws.Range("M2").Value = -508.8276
$ws.Range("H5").Value = 740.2
$ws.Range("I5").Value = 800.25
$ws.Range("K5").Value = 800.25
$ws.Range("M5").Value = -688.25
$ws.Range("H32").Value = 16104.377
$ws.Range("I32").Value = 15730.5
$ws.Range("J32").Value = 23332.666
$ws.Range("K32").Value = 15730.5
$ws.Range("L32").Value = 23332.666
$ws.Range("M32").Value = -15443.5
$ws.Range("N32").Value = -23906.666
$ws.Range("H61").Value = 5302.7354
$ws.Range("I61").Value = 5015.769
$ws.Range("K61").Value = 5015.769
$ws.Range("M61").Value = -4803.769
$ws.Range("H74").Value = 1597.4166
$ws.Range("I74").Value = 1207.125
$ws.Range("K74").Value = 1207.125
$ws.Range("M74").Value = -333.125
$ws.Range("H77").Value = 1597.4166
$ws.Range("I77").Value = 1207.125
$ws.Range("K77").Value = 6035.625
$ws.Range("M77").Value = -1667.625
$ws.Range("H116").Value = 700.9375
$ws.Range("I116").Value = 621.8276
$ws.Range("K116").Value = 621.8276
$ws.Range("M116").Value = 1672.1724
$ws.Range("H122").Value = 4258.143
$ws.Range("I122").Value = 3501.2727
$ws.Range("K122").Value = 10503.8181
$ws.Range("M122").Value = -8053.8181
$ws.Range("H132").Value = 19918.027
$ws.Range("I132").Value = 23553.809
$ws.Range("J132").Value = 11324.363
$ws.Range("K132").Value = 70661.427
$ws.Range("L132").Value = 33973.089
$ws.Range("M132").Value = -68131.427
$ws.Range("N132").Value = -39033.089
$ws.Range("H136").Value = 5302.7354
$ws.Range("I136").Value = 5015.769
$ws.Range("K136").Value = 15047.307
$ws.Range("M136").Value = -12497.307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 700.9375
$ws.Range("I3").Value = 621.8276
$ws.Range("K3").Value = 621.8276
$ws.Range("M3").Value = -507.8276
$ws.Range("H4").Value = 740.2
$ws.Range("I4").Value = 800.25
$ws.Range("K4").Value = 800.25
$ws.Range("M4").Value = -685.25
$ws.Range("H26").Value = 7379.4443
$ws.Range("I26").Value = 7379.4443
$ws.Range("K26").Value = 7379.4443
$ws.Range("M26").Value = -7087.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 139.75
$ws.Range("I7").Value = 153.77777
$ws.Range("J7").Value = 128.27272
$ws.Range("K7").Value = 153.77777
$ws.Range("L7").Value = 128.27272
$ws.Range("M7").Value = -40.77777
$ws.Range("N7").Value = -354.27272
$ws.Range("H22").Value = 1124.4286
$ws.Range("I22").Value = 1394.2
$ws.Range("K22").Value = 1394.2
$ws.Range("M22").Value = -1044.2
$ws.Range("H62").Value = 39139.727
$ws.Range("I62").Value = 19754.625
$ws.Range("J62").Value = 90833.336
$ws.Range("K62").Value = 19754.625
$ws.Range("L62").Value = 90833.336
$ws.Range("M62").Value = -19130.625
$ws.Range("N62").Value = -92081.336
$ws.Range("H65").Value = 39139.727
$ws.Range("I65").Value = 19754.625
$ws.Range("J65").Value = 90833.336
$ws.Range("K65").Value = 98773.125
$ws.Range("L65").Value = 454166.68
$ws.Range("M65").Value = -95653.125
$ws.Range("N65").Value = -460406.68
$ws.Range("H86").Value = 4665.9653
$ws.Range("I86").Value = 4251.2104
$ws.Range("J86").Value = 5454
$ws.Range("K86").Value = 4251.2104
$ws.Range("L86").Value = 5454
$ws.Range("M86").Value = -3128.2104
$ws.Range("N86").Value = -7700
$ws.Range("H89").Value = 4665.9653
$ws.Range("I89").Value = 4251.2104
$ws.Range("J89").Value = 5454
$ws.Range("K89").Value = 21256.052
$ws.Range("L89").Value = 27270
$ws.Range("M89").Value = -15640.052
$ws.Range("N89").Value = -38502
$ws.Range("H141").Value = 190848
$ws.Range("J141").Value = 202329.06
$ws.Range("L141").Value = 202329.06
$ws.Range("N141").Value = -212689.06

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 87.40000000000001
$ws.Range("I2").Value = 38
$ws.Range("J2").Value = 108.57143
$ws.Range("K2").Value = 228
$ws.Range("L2").Value = 651.42858
$ws.Range("M2").Value = -115
$ws.Range("N2").Value = -877.42858
$ws.Range("H8").Value = 251.25
$ws.Range("I8").Value = 251.25
$ws.Range("K8").Value = 753.75
$ws.Range("M8").Value = -614.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 20500
$ws.Range("J57").Value = 20500
$ws.Range("L57").Value = 20500
$ws.Range("N57").Value = -22140
$ws.Range("H102").Value = 530334.75
$ws.Range("I102").Value = 623315.2
$ws.Range("K102").Value = 623315.2
$ws.Range("M102").Value = -621693.2
$ws.Range("H126").Value = 1811.6897
$ws.Range("I126").Value = 1574.6538
$ws.Range("K126").Value = 4723.9614
$ws.Range("M126").Value = -2253.9614
$ws.Range("H141").Value = 40476
$ws.Range("J141").Value = 40476
$ws.Range("L141").Value = 40476
$ws.Range("N141").Value = -50836

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7823.6
$ws.Range("I7").Value = 7184
$ws.Range("K7").Value = 7184
$ws.Range("M7").Value = -7072
$ws.Range("H40").Value = 3949.5
$ws.Range("I40").Value = 3949.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3949.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3813.5
$ws.Range("N40").ClearContents() | Out-Null
$ws.Range("H46").Value = 4478.2104
$ws.Range("I46").Value = 5498.5
$ws.Range("J46").Value = 4358.1763
$ws.Range("K46").Value = 5498.5
$ws.Range("L46").Value = 4358.1763
$ws.Range("M46").Value = -5310.5
$ws.Range("N46").Value = -4734.1763
$ws.Range("H56").Value = 27504
$ws.Range("I56").Value = 6250
$ws.Range("K56").Value = 6250
$ws.Range("M56").Value = -5559
$ws.Range("H61").Value = 3139.3125
$ws.Range("I61").Value = 2944.5715
$ws.Range("K61").Value = 2944.5715
$ws.Range("M61").Value = -2742.5715
$ws.Range("H113").Value = 3139.3125
$ws.Range("I113").Value = 2944.5715
$ws.Range("K113").Value = 2944.5715
$ws.Range("M113").Value = -774.5715
$ws.Range("H126").Value = 7823.6
$ws.Range("I126").Value = 7184
$ws.Range("K126").Value = 21552
$ws.Range("M126").Value = -19082
$ws.Range("H140").Value = 76926.7
$ws.Range("J140").Value = 76926.7
$ws.Range("L140").Value = 76926.7
$ws.Range("N140").Value = -87286.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 35000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents() | Out-Null
$ws.Range("H132").Value = 10104323
$ws.Range("I132").Value = 15152491
$ws.Range("J132").Value = 7984.727
$ws.Range("K132").Value = 45457473
$ws.Range("L132").Value = 23954.181
$ws.Range("M132").Value = -45454943
$ws.Range("N132").Value = -29014.181
$ws.Range("H136").Value = 4389.425
$ws.Range("I136").Value = 2389.6155
$ws.Range("K136").Value = 7168.8465
$ws.Range("M136").Value = -4618.8465
$ws.Range("H138").Value = 98428.5
$ws.Range("J138").Value = 98428.5
$ws.Range("L138").Value = 98428.5
$ws.Range("N138").Value = -108708.5
$ws.Range("H140").Value = 98209.5
$ws.Range("J140").Value = 98209.5
$ws.Range("L140").Value = 98209.5
$ws.Range("N140").Value = -108569.5
